$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine last used row in column C (the "Förändrad" date column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

# Update every row's "Förändrad" date value by one day (45177 -> 45178),
# mirroring an automatic "last changed" timestamp bump applied to the whole column.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $cell.Value2 + 1
    }
}
